# Update the Xt-EHR "Procedure" mappable sheet to reflect the latest
# EHDSProcedure model definitions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: EHDSProcedure.identifier ---
$ws.Range("C3").Value = "Identifier for the procedure"
$ws.Range("D3").Value = "Identifier for the procedure"

# --- Row 4: becomes EHDSProcedure.subject (was .description) ---
$ws.Range("A4").Value = "EHDSProcedure.subject"
$ws.Range("B4").Value = "EHDSProcedure.subject"
$ws.Range("C4").Value = "On whom the procedure was performed."
$ws.Range("D4").Value = "On whom the procedure was performed."
$ws.Range("E4").Value = "EHDSPatient"
$ws.Range("F4").Value = "1..1"

# --- Row 5: EHDSProcedure.code ---
$ws.Range("C5").Value = "Code identifying the procedure"

# --- Row 6: EHDSProcedure.date[x] ---
$ws.Range("C6").Value = "Date and time of the procedure or interval of its performance"

# --- Row 7: EHDSProcedure.performer ---
$ws.Range("C7").Value = "An actor who performed the procedure"
$ws.Range("D7").Value = "An actor who performed the procedure"
$ws.Range("E7").Value = "EHDSHealthProfessional"

# --- Row 8: becomes EHDSProcedure.bodySite (was .anatomicLocation) ---
$ws.Range("A8").Value = "EHDSProcedure.bodySite"
$ws.Range("B8").Value = "EHDSProcedure.bodySite"
$ws.Range("C8").Value = "Anatomic location and laterality where the procedure was performed. This is the target site."
$ws.Range("E8").Value = "EHDSBodyStructure"
$ws.Range("G8").ClearContents()

# --- Row 9: becomes EHDSProcedure.reason[x] (was .reason) ---
$ws.Range("A9").Value = "EHDSProcedure.reason[x]"
$ws.Range("B9").Value = "EHDSProcedure.reason[x]"
$ws.Range("C9").Value = "The reason why the procedure was performed."
$ws.Range("E9").Value = "EHDSCondition"

# --- Row 10: EHDSProcedure.outcome ---
$ws.Range("C10").Value = "The outcome of the procedure - did it resolve the reasons for the procedure being performed?"

# --- Row 11: EHDSProcedure.complication ---
$ws.Range("C11").Value = "Any complications that occurred during the procedure, or in the immediate post-performance period. These are generally tracked separately from the procedure description, which will typically describe the procedure itself rather than any 'post procedure' issues."

# --- Row 12: EHDSProcedure.deviceUsed ---
$ws.Range("C12").Value = "Device used to perform the procedure"

# --- Row 13: EHDSProcedure.focalDevice ---
$ws.Range("C13").Value = "Device(s) that is/are implanted, removed, or otherwise manipulated (calibration, battery replacement, fitting a prosthesis, attaching a wound-vac, etc.) as a focal portion of the Procedure."
$ws.Range("D13").Value = "Device(s) that is/are implanted, removed, or otherwise manipulated (calibration, battery replacement, fitting a prosthesis, attaching a wound-vac, etc.) as a focal portion of the Procedure."
$ws.Range("E13").Value = "EHDSDevice"

# --- Row 14: EHDSProcedure.location ---
$ws.Range("C14").Value = "Location where the procedure was performed"
$ws.Range("D14").Value = "Location where the procedure was performed"

# --- Row 15: EHDSProcedure.note ---
$ws.Range("C15").Value = "Additional information about the procedure"

# --- Row 16 (old EHDSProcedure.subject) is removed: its content is now
#     merged into row 4 above. ---
$ws.Rows(16).Delete()
